$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 12501862
$ws.Range("J17").Value = 12501862
$ws.Range("L17").Value = 37505586
$ws.Range("N17").Value = -37505922
$ws.Range("H62").Value = 7816.909
$ws.Range("J62").Value = 3599.6
$ws.Range("L62").Value = 3599.6
$ws.Range("N62").Value = -4847.6
$ws.Range("H65").Value = 7816.909
$ws.Range("J65").Value = 3599.6
$ws.Range("L65").Value = 17998
$ws.Range("N65").Value = -24238
$ws.Range("H129").Value = 1295.25
$ws.Range("I129").Value = 843.125
$ws.Range("K129").Value = 2529.375
$ws.Range("M129").Value = 2470.625
$ws.Range("H137").Value = 3999.75
$ws.Range("I137").Value = 2666.3333
$ws.Range("J137").Value = 8000
$ws.Range("K137").Value = 7998.999899999999
$ws.Range("L137").Value = 24000
$ws.Range("M137").Value = -5448.999899999999
$ws.Range("N137").Value = -29100
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 275.83334
$ws.Range("I5").Value = 71.3
$ws.Range("J5").Value = 531.5
$ws.Range("K5").Value = 71.3
$ws.Range("L5").Value = 531.5
$ws.Range("M5").Value = 40.7
$ws.Range("N5").Value = -755.5
$ws.Range("H45").Value = 109514.58
$ws.Range("I45").Value = 156752.47
$ws.Range("J45").Value = 7165.8335
$ws.Range("K45").Value = 156752.47
$ws.Range("L45").Value = 7165.8335
$ws.Range("M45").Value = -156375.47
$ws.Range("N45").Value = -7919.8335
$ws.Range("H61").Value = 6834.8125
$ws.Range("I61").Value = 7285.8887
$ws.Range("J61").Value = 4399
$ws.Range("K61").Value = 7285.8887
$ws.Range("L61").Value = 4399
$ws.Range("M61").Value = -7073.8887
$ws.Range("N61").Value = -4823
$ws.Range("H74").Value = 7416.643
$ws.Range("I74").Value = 4426
$ws.Range("K74").Value = 4426
$ws.Range("M74").Value = -3552
$ws.Range("H77").Value = 7416.643
$ws.Range("I77").Value = 4426
$ws.Range("K77").Value = 22130
$ws.Range("M77").Value = -17762
$ws.Range("H97").Value = 1256.8182
$ws.Range("I97").Value = 1034.3529
$ws.Range("K97").Value = 1034.3529
$ws.Range("M97").Value = -538.3529000000001
$ws.Range("H110").Value = 4708.355
$ws.Range("I110").Value = 4166.0415
$ws.Range("J110").Value = 6567.7144
$ws.Range("K110").Value = 4166.0415
$ws.Range("L110").Value = 6567.7144
$ws.Range("M110").Value = -2121.0415
$ws.Range("N110").Value = -10657.7144
$ws.Range("H136").Value = 6834.8125
$ws.Range("I136").Value = 7285.8887
$ws.Range("J136").Value = 4399
$ws.Range("K136").Value = 21857.6661
$ws.Range("L136").Value = 13197
$ws.Range("M136").Value = -19307.6661
$ws.Range("N136").Value = -18297
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 275.83334
$ws.Range("I4").Value = 71.3
$ws.Range("J4").Value = 531.5
$ws.Range("K4").Value = 71.3
$ws.Range("L4").Value = 531.5
$ws.Range("M4").Value = 43.7
$ws.Range("N4").Value = -761.5
$ws.Range("H7").Value = 66666
$ws.Range("I7").Value = 66666
$ws.Range("K7").Value = 66666
$ws.Range("M7").Value = -66553
$ws.Range("H99").Value = 3871.2942
$ws.Range("I99").Value = 3053.4783
$ws.Range("K99").Value = 3053.4783
$ws.Range("M99").Value = -1555.4783
$ws.Range("H105").Value = 1164.0834
$ws.Range("I105").Value = 997.2727
$ws.Range("K105").Value = 997.2727
$ws.Range("M105").Value = 749.7273
$ws.Range("H107").Value = 1282.8948
$ws.Range("I107").Value = 1429.2727
$ws.Range("K107").Value = 1429.2727
$ws.Range("M107").Value = 490.7273
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 74.40000000000001
$ws.Range("J7").Value = 69.666664
$ws.Range("L7").Value = 69.666664
$ws.Range("N7").Value = -295.666664
$ws.Range("H31").Value = 6071.407
$ws.Range("I31").Value = 5173.5
$ws.Range("J31").Value = 6449.4736
$ws.Range("K31").Value = 5173.5
$ws.Range("L31").Value = 6449.4736
$ws.Range("M31").Value = -4878.5
$ws.Range("N31").Value = -7039.4736
$ws.Range("H34").Value = 6071.407
$ws.Range("I34").Value = 5173.5
$ws.Range("J34").Value = 6449.4736
$ws.Range("K34").Value = 5173.5
$ws.Range("L34").Value = 6449.4736
$ws.Range("M34").Value = -4971.5
$ws.Range("N34").Value = -6853.4736
$ws.Range("H58").Value = 12601.75
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H99").Value = 8690
$ws.Range("I99").Value = 9028
$ws.Range("J99").Value = 8216.799999999999
$ws.Range("K99").Value = 9028
$ws.Range("L99").Value = 8216.799999999999
$ws.Range("M99").Value = -7530
$ws.Range("N99").Value = -11212.8
$ws.Range("H122").Value = 3097.2222
$ws.Range("I122").Value = 3177.0527
$ws.Range("K122").Value = 9531.158100000001
$ws.Range("M122").Value = -7081.158100000001
$ws.Range("H126").Value = 8690
$ws.Range("I126").Value = 9028
$ws.Range("J126").Value = 8216.799999999999
$ws.Range("K126").Value = 27084
$ws.Range("L126").Value = 24650.4
$ws.Range("M126").Value = -24614
$ws.Range("N126").Value = -29590.4
$ws.Range("H132").Value = 4483.421
$ws.Range("I132").Value = 2574.125
$ws.Range("K132").Value = 7722.375
$ws.Range("M132").Value = -5192.375
$ws.Range("H135").Value = 90122.664
$ws.Range("J135").Value = 90122.664
$ws.Range("L135").Value = 90122.664
$ws.Range("N135").Value = -100262.664
$ws.Range("H136").Value = 12601.75
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("H141").Value = 257146.62
$ws.Range("J141").Value = 257146.62
$ws.Range("L141").Value = 257146.62
$ws.Range("N141").Value = -267506.62
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 122.85714
$ws.Range("I61").Value = 65
$ws.Range("K61").Value = 195
$ws.Range("M61").Value = 20
$ws.Range("H122").Value = 1198.3334
$ws.Range("I122").Value = 1400
$ws.Range("K122").Value = 12600
$ws.Range("M122").Value = -10150
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2665.8333
$ws.Range("I80").Value = 2374
$ws.Range("J80").Value = 3249.5
$ws.Range("K80").Value = 2374
$ws.Range("L80").Value = 3249.5
$ws.Range("M80").Value = -1376
$ws.Range("N80").Value = -5245.5
$ws.Range("H83").Value = 2665.8333
$ws.Range("I83").Value = 2374
$ws.Range("J83").Value = 3249.5
$ws.Range("K83").Value = 11870
$ws.Range("L83").Value = 16247.5
$ws.Range("M83").Value = -6878
$ws.Range("N83").Value = -26231.5
$ws.Range("H102").Value = 5191.731
$ws.Range("I102").Value = 2948.8235
$ws.Range("K102").Value = 2948.8235
$ws.Range("M102").Value = -1326.8235
$ws.Range("H107").Value = 835.0714
$ws.Range("I107").Value = 460.5
$ws.Range("K107").Value = 460.5
$ws.Range("M107").Value = 1459.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4759.7856
$ws.Range("I7").Value = 4672.5454
$ws.Range("K7").Value = 4672.5454
$ws.Range("M7").Value = -4560.5454
$ws.Range("H22").Value = 1971
$ws.Range("I22").Value = 2406.8572
$ws.Range("J22").Value = 1462.5
$ws.Range("K22").Value = 2406.8572
$ws.Range("L22").Value = 1462.5
$ws.Range("M22").Value = -2111.8572
$ws.Range("N22").Value = -2052.5
$ws.Range("H27").Value = 1971
$ws.Range("I27").Value = 2406.8572
$ws.Range("J27").Value = 1462.5
$ws.Range("K27").Value = 2406.8572
$ws.Range("L27").Value = 1462.5
$ws.Range("M27").Value = -2299.8572
$ws.Range("N27").Value = -1676.5
$ws.Range("H46").Value = 5089.8
$ws.Range("I46").Value = 1816.3334
$ws.Range("K46").Value = 1816.3334
$ws.Range("M46").Value = -1628.3334
$ws.Range("H55").Value = 45454652
$ws.Range("J55").Value = 133.25
$ws.Range("L55").Value = 133.25
$ws.Range("N55").Value = -479.25
$ws.Range("H61").Value = 65004.625
$ws.Range("I61").Value = 79660
$ws.Range("K61").Value = 79660
$ws.Range("M61").Value = -79458
$ws.Range("H113").Value = 65004.625
$ws.Range("I113").Value = 79660
$ws.Range("K113").Value = 79660
$ws.Range("M113").Value = -77490
$ws.Range("H122").Value = 2240.6365
$ws.Range("I122").Value = 2240.6365
$ws.Range("K122").Value = 6721.9095
$ws.Range("M122").Value = -4271.9095
$ws.Range("H126").Value = 4759.7856
$ws.Range("I126").Value = 4672.5454
$ws.Range("K126").Value = 14017.6362
$ws.Range("M126").Value = -11547.6362
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6334
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 6334
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 50000
$ws.Range("N65").Value = -56240

Write-Output "applied"